$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set headers for new columns, copying style from an existing header cell
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 29).Value = 85
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 0
}
